$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Build the new query strings (exact text, LF line endings, no trailing newline) ---

$casesQuery = "MATCH (c:case)`n MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)`n MATCH (f:file)-[*]->(c)`nWHERE f.file_type = 'Aligned RNA reads file' `nRETURN DISTINCT`n    c.case_id AS ``Case ID``,`n     ct.clinical_trial_designation AS ``Trial Code``,`n     a.arm_id AS Arm,`n      a.arm_drug AS ``Arm Treatment``,`nc.disease AS Diagnosis,`n  c.gender AS Gender,`n    c.race AS Race,`n    c.ethnicity AS Ethnicity"

$statQuery = "MATCH (f:file)`nOPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)`nOPTIONAL MATCH (f)-[*]->(c:case)`nWITH f,a,ct,c`n      WHERE f.file_type = 'Aligned RNA reads file'`nRETURN`n    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,`n    COUNT(DISTINCT c.case_id) AS Cases,`n    COUNT(DISTINCT f) AS Files"

$filesQuery = "MATCH (f:file)`nOPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)`nOPTIONAL MATCH (f)-[*]->(c:case)`nOPTIONAL MATCH (f)-->(parent)`nWITH f,a,ct,c,parent`nWHERE f.file_type = 'Aligned RNA reads file'`n WITH`n    f, parent, c, a, ct,`n    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n    toInteger(floor(log(f.file_size)/log(1024))) as i,`n    2 as precision`nWITH`n    f, parent, c, a, ct,`n    f.file_size /(1024^i) AS value,`n    10^precision AS factor,`n    units[i] as unit`nWITH`n    f, parent, c, a, ct, unit,`n    round(factor * value)/factor AS size`nRETURN DISTINCT`n    f.file_name AS ``File Name``,`n    head(labels(parent)) as Association,`n    f.file_description AS Description,`n    f.file_format AS ``File Format``,`n    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n    ct.clinical_trial_designation AS ``Trial Code``,`n    a.arm_id AS Arm,`n    c.case_id AS ``Case ID``"

# --- Insert a new column before column A, shifting existing data to B:E ---
$ws.Range("A1").EntireColumn.Insert()

# --- Row 1 headers: new column A gets "TabName" ---
$ws.Range("A1").Value = "TabName"

# --- Row 2 becomes the "CasesTab" row ---
$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statQuery

# --- Row 3 is a new row: the "FilesTab" row ---
$ws.Range("A3").Value = "FilesTab"
$ws.Range("B3").Value = $filesQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("D3").Value = "TC02_Trials_Filter_AssocFileType-AlignedRNA_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC02_Trials_Filter_AssocFileType-AlignedRNA_WebData.xlsx"

# --- Apply the wrap-text style (same as used on B2) to the new query cells ---
$ws.Range("C2").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# --- Column widths ---
$ws.Range("A1").ColumnWidth = 7.983072916666667
$ws.Range("B1:C1").ColumnWidth = 74.98307291666667
$ws.Range("D1").ColumnWidth = 69.43619791666667
$ws.Range("E1").ColumnWidth = 27.709635416666668

# --- Row heights ---
$ws.Range("A2").RowHeight = 188.5
$ws.Range("A3").RowHeight = 409.5

# --- View state: scroll/selection ---
$ws.Range("C3").Select()
$excel.ActiveWindow.ScrollRow = 3
